$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'71.027.97"
$ws.Range("E2").Value = "  -0.05%  "

# Row 3
$ws.Range("D3").Value = "'3.832.62"
$ws.Range("E3").Value = "  +0.50%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'708.50"
$ws.Range("E5").Value = "  +1.26%  "

# Row 6
$ws.Range("D6").Value = "'172.79"
$ws.Range("E6").Value = "  -0.54%  "

# Row 7
$ws.Range("D7").Value = "'3.829.62"
$ws.Range("E7").Value = "  +0.37%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("E9").Value = "  -0.32%  "

# Row 10
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  -0.21%  "

# Row 11
$ws.Range("D11").Value = "'7.33"
$ws.Range("E11").Value = "  +1.51%  "

# Row 12
$ws.Range("E12").Value = "  -0.17%  "

# Row 13
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  -1.11%  "

# Row 14
$ws.Range("D14").Value = "'36.67"
$ws.Range("E14").Value = "  +1.01%  "

# Row 15
$ws.Range("D15").Value = "'4.478.73"
$ws.Range("E15").Value = "  +0.56%  "

# Row 16
$ws.Range("D16").Value = "'3.813.73"
$ws.Range("E16").Value = "  +0.31%  "

# Row 17
$ws.Range("D17").Value = "'71.016.00"
$ws.Range("E17").Value = "  +0.01%  "

# Row 18
$ws.Range("D18").Value = "'7.21"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19
$ws.Range("E19").Value = "  +0.60%  "

# Row 20
$ws.Range("D20").Value = "'17.34"
$ws.Range("E20").Value = "  -2.71%  "

# Row 21
$ws.Range("D21").Value = "'10.76"
$ws.Range("E21").Value = "  -3.15%  "

# Row 22
$ws.Range("D22").Value = "'494.14"
$ws.Range("E22").Value = "  +3.00%  "

# Row 23
$ws.Range("D23").Value = "'0.721"
$ws.Range("E23").Value = "  +1.33%  "

# Row 24
$ws.Range("D24").Value = "'84.99"
$ws.Range("E24").Value = "  +1.36%  "

# Row 25
$ws.Range("D25").Value = "'0.0000146"
$ws.Range("E25").Value = "  +1.42%  "

# Row 26
$ws.Range("D26").Value = "'10.62"
$ws.Range("E26").Value = "  +1.63%  "

# Row 27
$ws.Range("D27").Value = "'12.16"
$ws.Range("E27").Value = "  -1.35%  "

# Row 28
$ws.Range("D28").Value = "'2.10"
$ws.Range("E28").Value = "  -3.10%  "

# Row 29
$ws.Range("E29").Value = "  +2.20%  "

# Row 30
$ws.Range("E30").Value = "  -0.04%  "

# Row 31
$ws.Range("D31").Value = "'7.50"
$ws.Range("E31").Value = "  -0.40%  "

# Row 32
$ws.Range("D32").Value = "'2.25"
$ws.Range("E32").Value = "  -2.23%  "

# Row 33
$ws.Range("D33").Value = "'29.40"
$ws.Range("E33").Value = "  -0.45%  "

# Row 34
$ws.Range("D34").Value = "'0.178"
$ws.Range("E34").Value = "  -4.60%  "

# Row 35
$ws.Range("D35").Value = "'9.19"
$ws.Range("E35").Value = "  -0.77%  "

# Row 36
$ws.Range("D36").Value = "'3.798.08"
$ws.Range("E36").Value = "  +0.93%  "

# Row 37
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("E38").Value = "  -0.29%  "

# Row 39
$ws.Range("D39").Value = "'2.33"
$ws.Range("E39").Value = "  +3.62%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'6.01"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "'1.03"
$ws.Range("E41").Value = "  +5.23%  "

# Row 42
$ws.Range("D42").Value = "'3.35"
$ws.Range("E42").Value = "  -2.41%  "

# Row 43
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("E44").Value = "  +0.15%  "

# Row 45
$ws.Range("D45").Value = "'0.000314"
$ws.Range("E45").Value = "  -3.97%  "

# Row 46
$ws.Range("D46").Value = "'163.35"
$ws.Range("E46").Value = "  +0.12%  "

# Row 47
$ws.Range("D47").Value = "'48.70"
$ws.Range("E47").Value = "  -0.52%  "

# Row 48
$ws.Range("E48").Value = "  +0.51%  "

# Row 49
$ws.Range("D49").Value = "'414.69"
$ws.Range("E49").Value = "  +1.18%  "

# Row 50
$ws.Range("E50").Value = "  +0.61%  "

# Row 51
$ws.Range("D51").Value = "'0.296"
$ws.Range("E51").Value = "  -1.75%  "
